$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.669.27'
$ws.Range('E2').Value = '  +1.23%  '
$ws.Range('D3').Value = '1.827.65'
$ws.Range('E3').Value = '  +1.85%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.008'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.007'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '308.01'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4645'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.50%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3608'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07136'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9046'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07769'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.42'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').Value = '1.850.41'
$ws.Range('E13').Value = '  +2.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.268'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.44%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.349'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '87.72'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.37%  '
$ws.Range('E17').Value = '  +0.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008569'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.006'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.12%  '
$ws.Range('D20').Value = '26.720.81'
$ws.Range('E20').Value = '  +1.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.19'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.53%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.018'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.59%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.56'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.927'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.65'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('E26').Value = '  +0.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.977'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.59%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '114.14'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.84%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.826'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.94%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08811'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.148'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.59%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7337'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.59%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.143'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.25%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.443'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.25%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.728'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.077'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.61%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01926'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.929'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.82%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05137'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.63%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.877'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5064'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1498'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.027'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.44%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4674'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.94%  '
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.04'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.73%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '98.42'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.63%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.563'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06034'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '64.02'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '35.92'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.44%  '
